# Fill in the missing time-log entry on row 28 (date 2017-11-20, a second
# shift that day) and move the view/selection the way the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A28: Date, B28: Start Time, C28: End Time (all stored as Excel serials,
# matching how the rest of the log's rows are populated).
$ws.Range("A28").Value = 43059
$ws.Range("B28").Value = 0.79513888888888884
$ws.Range("C28").Value = 0.87847222222222221

# D28 already carries the shared formula =ABS(C28-B28); force a recalc so
# its cached value (and the SUM() total in D36) reflect the new entry.
$excel.Calculate()

# Move the viewport/selection to match where the author ended up editing.
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D24:D28").Select()
